$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.226.14'
$ws.Range("E2").Value = '  -0.63%  '

# Row 3
$ws.Range("D3").Value = '2.354.18'
$ws.Range("E3").Value = '  +4.79%  '

# Row 4
$ws.Range("E4").Value = '  +0.22%  '

# Row 5
$ws.Range("D5").Value = '''234.34'
$ws.Range("E5").Value = '  +1.79%  '

# Row 6
$ws.Range("E6").Value = '  +1.35%  '

# Row 7
$ws.Range("D7").Value = '''71.73'
$ws.Range("E7").Value = '  +13.73%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '''0.495'
$ws.Range("E9").Value = '  +12.15%  '

# Row 10
$ws.Range("D10").Value = '''0.0974'
$ws.Range("E10").Value = '  +1.64%  '

# Row 11
$ws.Range("D11").Value = '''27.38'
$ws.Range("E11").Value = '  +0.30%  '

# Row 12
$ws.Range("E12").Value = '  +2.33%  '

# Row 13
$ws.Range("D13").Value = '2.706.64'
$ws.Range("E13").Value = '  +5.18%  '

# Row 14
$ws.Range("E14").Value = '  +5.45%  '

# Row 15
$ws.Range("D15").Value = '''6.36'
$ws.Range("E15").Value = '  +4.98%  '

# Row 16
$ws.Range("D16").Value = '''0.866'
$ws.Range("E16").Value = '  +4.74%  '

# Row 17
$ws.Range("D17").Value = '2.350.16'
$ws.Range("E17").Value = '  +5.31%  '

# Row 18
$ws.Range("D18").Value = '43.264.06'
$ws.Range("E18").Value = '  -0.02%  '

# Row 19
$ws.Range("D19").Value = '''0.0000100'
$ws.Range("E19").Value = '  +4.29%  '

# Row 20
$ws.Range("E20").Value = '  +4.40%  '

# Row 21
$ws.Range("D21").Value = '''74.65'
$ws.Range("E21").Value = '  +2.44%  '

# Row 22
$ws.Range("D22").Value = '''250.10'
$ws.Range("E22").Value = '  +1.50%  '

# Row 23
$ws.Range("D23").Value = '''3.83'
$ws.Range("E23").Value = '  +3.86%  '

# Row 24
$ws.Range("E24").Value = '  -0.02%  '

# Row 25
$ws.Range("E25").Value = '  +2.03%  '

# Row 26
$ws.Range("D26").Value = '''2.34'
$ws.Range("E26").Value = '  +2.68%  '

# Row 27
$ws.Range("E27").Value = '  +3.42%  '

# Row 28
$ws.Range("D28").Value = '''22.41'
$ws.Range("E28").Value = '  +4.09%  '

# Row 29
$ws.Range("D29").Value = '''172.69'
$ws.Range("E29").Value = '  +0.37%  '

# Row 30
$ws.Range("E30").Value = '  +9.18%  '

# Row 31
$ws.Range("E31").Value = '  +1.36%  '

# Row 32
$ws.Range("E32").Value = '  +2.68%  '

# Row 33
$ws.Range("D33").Value = '''5.01'
$ws.Range("E33").Value = '  +2.75%  '

# Row 34
$ws.Range("E34").Value = '  +2.74%  '

# Row 35
$ws.Range("D35").Value = '''5.08'
$ws.Range("E35").Value = '  +4.52%  '

# Row 36
$ws.Range("E36").Value = '  +3.21%  '

# Row 37
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").Value = '''2.43'
$ws.Range("E37").Value = '  +7.08%  '

# Row 38
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '''6.54'
$ws.Range("E38").Value = '  +4.19%  '

# Row 39
$ws.Range("E39").Value = '  +1.83%  '

# Row 40
$ws.Range("D40").Value = '''19.00'
$ws.Range("E40").Value = '  +12.41%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''8.94'
$ws.Range("E41").Value = '  +4.29%  '

# Row 42
$ws.Range("B42").Value = 'BinanceUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("E43").Value = '  -0.23%  '

# Row 44
$ws.Range("D44").Value = '''99.34'
$ws.Range("E44").Value = '  +3.06%  '

# Row 45
$ws.Range("E45").Value = '  +9.27%  '

# Row 46
$ws.Range("D46").Value = '''0.0960'
$ws.Range("E46").Value = '  +2.15%  '

# Row 47
$ws.Range("E47").Value = '  +3.06%  '

# Row 48
$ws.Range("D48").Value = '1.445.79'
$ws.Range("E48").Value = '  -0.17%  '

# Row 49
$ws.Range("D49").Value = '2.578.92'
$ws.Range("E49").Value = '  +5.30%  '

# Row 50
$ws.Range("D50").Value = '''2.76'
$ws.Range("E50").Value = '  +0.85%  '

# Row 51
$ws.Range("E51").Value = '  -4.57%  '
